$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Part 1: Replace the body text of the "What I Learned" paragraph
#         (the paragraph that used to contain just "..." / ".") with the
#         full writeup text describing the CarQueue lab.
# -----------------------------------------------------------------------
$newText = "I began this lab by creating a class called CarQueue that would represent and maintain a queue of random directions that cars should travel. To implement the queue I used the ArrayDeque class and set the generic type parameter to be Integer. I knew from the instructions to create an addToQueue() method and a deleteQueue() method. I also knew that I should create a constructor method that would populate the queue with 6 random integers. The random integers that CarQueue can generate are 0 through 3, where 0 signifies traveling up, 1 for down, 2 for right, and 3 for left. "

$p3 = $d.Paragraphs(3).Range
$body = $d.Range($p3.Start, $p3.End - 1)
$body.Text = $newText

# -----------------------------------------------------------------------
# Part 2: Relocate the "_GoBack" bookmark. It used to wrap the paragraph
#         holding the GitHub screenshot image; it now sits as a
#         zero-length bookmark right after "...for left." and before the
#         trailing space that closes the paragraph we just rewrote.
# -----------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$p3after = $d.Paragraphs(3).Range
$markPos = $p3after.End - 2
$markRange = $d.Range($markPos, $markPos)
$d.Bookmarks.Add("_GoBack", $markRange)
